# Burndown chart update - "include progress from a member"
# Updates the ANTICIPATED REMAINING values (B3:B7) to a tapering burn-down
# and lowers ACTUAL REMAINING (C7:C8) to reflect the extra progress,
# then reselects the anticipated-remaining column and switches the sheet
# to portrait print orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ANTICIPATED REMAINING (column B) - tapering values instead of flat 25s
$ws.Cells.Item(3, 2).Value = 20.833333333333329
$ws.Cells.Item(4, 2).Value = 16.666666666666661
$ws.Cells.Item(5, 2).Value = 12.499999999999991
$ws.Cells.Item(6, 2).Value = 8.3333333333333215
$ws.Cells.Item(7, 2).Value = 4.1666666666666501

# ACTUAL REMAINING (column C) - member progress drops remaining from 19 to 17
$ws.Cells.Item(7, 3).Value = 17
$ws.Cells.Item(8, 3).Value = 17

# Page setup: print orientation -> portrait
$ws.PageSetup.Orientation = 1

# Update the active selection to the full anticipated-remaining range
$ws.Range("B2:B8").Select()
